# add synthesized alu with abc and update results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Apply "#,##0" number formatting to the Synopsys-PCF block (P:U) for all
#    data rows. This reuses the existing styles (s=7 for P:T, s=6 for U,
#    which already carry the bold/plain fonts) instead of minting new ones.
# ---------------------------------------------------------------------------
$ws.Range("P5:T6").NumberFormat = "#,##0"
$ws.Range("U5:U6").NumberFormat = "#,##0"

$ws.Range("P8:T9").NumberFormat = "#,##0"
$ws.Range("U8:U9").NumberFormat = "#,##0"

$ws.Range("P12:T13").NumberFormat = "#,##0"
$ws.Range("U12:U13").NumberFormat = "#,##0"

$ws.Range("P15:T16").NumberFormat = "#,##0"
$ws.Range("U15:U16").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# 2. Extend the formatted-but-empty separator/blank rows (7, 10, 11, 14) so
#    that the P:U block carries the same formatting as the data rows above.
# ---------------------------------------------------------------------------
$ws.Range("P7:T7").NumberFormat = "#,##0"
$ws.Range("U7").NumberFormat = "#,##0"

$ws.Range("P10:T10").NumberFormat = "#,##0"
$ws.Range("U10").NumberFormat = "#,##0"

$ws.Range("P11:T11").NumberFormat = "#,##0"
$ws.Range("U11").NumberFormat = "#,##0"

$ws.Range("P14:T14").NumberFormat = "#,##0"
$ws.Range("U14").NumberFormat = "#,##0"

# Row 14 also grows a few more formatted-but-empty cells (mirrors the
# Yosys-ABC block used on the other data rows), plus plain (style 0) cells
# in A/B to match row 13's new blank A cell.
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("H14").NumberFormat = "#,##0"
$ws.Range("I14").NumberFormat = "#,##0"

$ws.Range("A13").Font.Bold = $false
$ws.Range("A14").Font.Bold = $false
$ws.Range("B14").Font.Bold = $false

# ---------------------------------------------------------------------------
# 3. New "alu" rows (15/16): fill in the synthesized ABC results (Yosys-ABC
#    block, columns D:I) that were previously blank.
# ---------------------------------------------------------------------------
$ws.Range("D15").Value = 41
$ws.Range("E15").Value = 111
$ws.Range("F15").Formula = "=D15+E15"
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 195
$ws.Range("H15").Formula = "=D15+E15+G15"
$ws.Range("H15").NumberFormat = "#,##0"
$ws.Range("I15").Formula = "=F15+(G15*5)"
$ws.Range("I15").NumberFormat = "#,##0"

$ws.Range("D16").Value = 211
$ws.Range("E16").Value = 346
$ws.Range("F16").Formula = "=D16+E16"
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("G16").Value = 685
$ws.Range("H16").Formula = "=D16+E16+G16"
$ws.Range("H16").NumberFormat = "#,##0"
$ws.Range("I16").Formula = "=F16+(G16*5)"
$ws.Range("I16").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# 4. View state: scroll back to the left edge and move the selection.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.TabRatio = 76
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("S25").Select()
